$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: *_old -> *_FV2210, *_new -> *_FV2304
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2210")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2304")
        }
    }
}

# Add table over the used range
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
